# Regenerate orders with updated distance/size codes.
# Distance codes embedded throughout Condition / Filename_Left / Filename_Right
# / Distance columns:   D64 -> D69,  D51 -> D55,  D80 -> D86
# Size code embedded the same way (S25 / S20 stay as-is): S30 -> S31
#
# All of these are substrings inside longer shared-string values (e.g.
# "Face09_D64_S30" / "Face09_D64_S30_l.png"), so the Find/Replace must match
# *part* of the cell text, not the whole cell. Range.Replace's 3rd
# positional arg is LookAt: 2 = xlPart (substring match), 1 = xlWhole.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$rng = $ws.UsedRange

$xlPart = 2

# None of the new codes (69/55/86/31) collide with any of the still-to-be-
# processed old codes (64/51/80/30), so one sequential pass is safe.
$rng.Replace("D64", "D69", $xlPart, 1, $false, $false, $false, $false)
$rng.Replace("D51", "D55", $xlPart, 1, $false, $false, $false, $false)
$rng.Replace("D80", "D86", $xlPart, 1, $false, $false, $false, $false)
$rng.Replace("S30", "S31", $xlPart, 1, $false, $false, $false, $false)
